# Applies the "Updated symbol list on Sun Dec 18 17:37:39 UTC 2022 with
# GitHub Actions" edit to the crypto price sheet.
#
# The sheet stores every data value (prices, labels, links, ids) as text
# (inline strings in the source OOXML). Many of the "Price" column values
# look like plain numbers (e.g. "247.58"), so writing them with a plain
# `.Value =` assignment would let Excel re-interpret them as numeric
# values and lose the original text formatting (trailing zeros, exact
# decimal representation, etc.). To keep them as text we force the
# cell's number format to Text ("@") before assigning the literal
# string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Price column (D) updates -------------------------------------------
Set-TextValue "D2"  "247.58"
Set-TextValue "D4"  "5.478"
Set-TextValue "D5"  "0.05621"
Set-TextValue "D6"  "6.460"
Set-TextValue "D7"  "0.8043"
Set-TextValue "D8"  "1.040"
Set-TextValue "D9"  "0.1423"
Set-TextValue "D10" "0.07295"
Set-TextValue "D12" "0.02921"
Set-TextValue "D13" "0.09252"
Set-TextValue "D14" "0.001676"
Set-TextValue "D15" "3.221"
Set-TextValue "D16" "0.04744"
Set-TextValue "D17" "0.0005830"
Set-TextValue "D18" "0.006443"
Set-TextValue "D19" "0.005069"
Set-TextValue "D23" "3.379"
Set-TextValue "D24" "2.122"

# Row 26 (ProBitToken) - only the "Volume(1h)" label cell changes.
$ws.Range("E26").Value = "25ProBitTokenPROBBestin24h"

Set-TextValue "D27" "0.0003305"
Set-TextValue "D40" "0.04172"
Set-TextValue "D41" "0.006903"

# --- Rows 42/43: CEJI and BKEXToken swap places --------------------------
# Row 42 becomes BKEXToken (previously CEJI)
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1039"
$ws.Range("E42").Value = "41BKEXTokenBKK"

# Row 43 becomes CEJI (previously BKEXToken)
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002975"
$ws.Range("E43").Value = "42CEJICEJI"

Set-TextValue "D44" "0.009000"
Set-TextValue "D45" "0.00005645"
Set-TextValue "D46" "0.00000000751"
Set-TextValue "D47" "0.6811"
Set-TextValue "D48" "0.01590"
Set-TextValue "D49" "0.00002104"
